$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 887.8570999999999
$ws.Range("I8").Value = 206
$ws.Range("J8").Value = 1399.25
$ws.Range("K8").Value = 618
$ws.Range("L8").Value = 4197.75
$ws.Range("M8").Value = -479
$ws.Range("N8").Value = -4475.75
$ws.Range("H31").Value = 669
$ws.Range("I31").Value = 669
$ws.Range("K31").Value = 2007
$ws.Range("M31").Value = -1777
$ws.Range("H32").Value = 10150
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H33").Value = 198
$ws.Range("I33").Value = 198
$ws.Range("K33").Value = 198
$ws.Range("M33").Value = 31
$ws.Range("H38").Value = 908.53845
$ws.Range("J38").Value = 1925
$ws.Range("L38").Value = 5775
$ws.Range("N38").Value = -6519
$ws.Range("H39").Value = 549.8889
$ws.Range("I39").Value = 16.4
$ws.Range("J39").Value = 1216.75
$ws.Range("K39").Value = 49.2
$ws.Range("L39").Value = 3650.25
$ws.Range("M39").Value = 246.8
$ws.Range("N39").Value = -4242.25
$ws.Range("H43").Value = 3338.4443
$ws.Range("I43").Value = 2674
$ws.Range("J43").Value = 3870
$ws.Range("K43").Value = 2674
$ws.Range("L43").Value = 3870
$ws.Range("M43").Value = -2605
$ws.Range("N43").Value = -4008
$ws.Range("H48").Value = 1125
$ws.Range("I48").Value = 500
$ws.Range("K48").Value = 1500
$ws.Range("M48").Value = -1208
$ws.Range("H51").Value = 6000
$ws.Range("I51").Value = 5000
$ws.Range("K51").Value = 5000
$ws.Range("M51").Value = -4516
$ws.Range("H55").Value = 279.8889
$ws.Range("J55").Value = 455.8
$ws.Range("L55").Value = 455.8
$ws.Range("N55").Value = -883.8
$ws.Range("H56").Value = 1125
$ws.Range("I56").Value = 500
$ws.Range("K56").Value = 1500
$ws.Range("M56").Value = -966
$ws.Range("H74").Value = 5488.8
$ws.Range("I74").Value = 4722
$ws.Range("K74").Value = 4722
$ws.Range("M74").Value = -3786
$ws.Range("H77").Value = 5488.8
$ws.Range("I77").Value = 4722
$ws.Range("K77").Value = 23610
$ws.Range("M77").Value = -18930
$ws.Range("H101").Value = 710.4
$ws.Range("I101").Value = 489
$ws.Range("J101").Value = 765.75
$ws.Range("K101").Value = 1467
$ws.Range("L101").Value = 2297.25
$ws.Range("M101").Value = 155
$ws.Range("N101").Value = -5541.25
$ws.Range("H132").Value = 4503.2354
$ws.Range("I132").Value = 3824.9285
$ws.Range("J132").Value = 7668.6665
$ws.Range("K132").Value = 11474.7855
$ws.Range("L132").Value = 23005.9995
$ws.Range("M132").Value = -8944.7855
$ws.Range("N132").Value = -28065.9995
$ws.Range("H135").Value = 249.5
$ws.Range("I135").Value = 249.5
$ws.Range("K135").Value = 2245.5
$ws.Range("M135").Value = 289.5
$ws.Range("H138").Value = 3267
$ws.Range("I138").Value = 2369.5
$ws.Range("K138").Value = 7108.5
$ws.Range("M138").Value = -1968.5

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5000
$ws.Range("I2").Value = 5000
$ws.Range("K2").Value = 5000
$ws.Range("M2").Value = -4887
$ws.Range("H32").Value = 2381.85
$ws.Range("I32").Value = 2485.5
$ws.Range("J32").Value = 1449
$ws.Range("K32").Value = 2485.5
$ws.Range("L32").Value = 1449
$ws.Range("M32").Value = -2198.5
$ws.Range("N32").Value = -2023
$ws.Range("H116").Value = 5000
$ws.Range("I116").Value = 5000
$ws.Range("K116").Value = 5000
$ws.Range("M116").Value = -2706
$ws.Range("H122").Value = 2654.4546
$ws.Range("I122").Value = 3066.3333
$ws.Range("K122").Value = 9198.999899999999
$ws.Range("M122").Value = -6748.999899999999
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("H135").Value = 99429
$ws.Range("J135").Value = 99429
$ws.Range("L135").Value = 99429
$ws.Range("N135").Value = -109569

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5000
$ws.Range("I3").Value = 5000
$ws.Range("K3").Value = 5000
$ws.Range("M3").Value = -4886
$ws.Range("H22").Value = 4863
$ws.Range("I22").Value = 4863
$ws.Range("K22").Value = 4863
$ws.Range("M22").Value = -4690
$ws.Range("H134").Value = 7042.5
$ws.Range("I134").Value = 6249.5
$ws.Range("K134").Value = 18748.5
$ws.Range("M134").Value = -16213.5

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 900
$ws.Range("I16").Value = 800
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 800
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -513
$ws.Range("N16").Value = -1574
$ws.Range("H41").Value = 20055
$ws.Range("J41").Value = 20055
$ws.Range("L41").Value = 20055
$ws.Range("N41").Value = -20911
$ws.Range("H50").Value = 28777.777
$ws.Range("J50").Value = 28428.572
$ws.Range("L50").Value = 28428.572
$ws.Range("N50").Value = -29678.572
$ws.Range("H51").Value = 24600
$ws.Range("J51").Value = 24600
$ws.Range("L51").Value = 24600
$ws.Range("N51").Value = -26072
$ws.Range("H59").Value = 30000
$ws.Range("I59").Value = 10000
$ws.Range("K59").Value = 10000
$ws.Range("M59").Value = -8855
$ws.Range("H60").Value = 24000
$ws.Range("I60").Value = 10000
$ws.Range("J60").Value = 26333.334
$ws.Range("K60").Value = 10000
$ws.Range("L60").Value = 26333.334
$ws.Range("M60").Value = -9489
$ws.Range("N60").Value = -27355.334
$ws.Range("H61").Value = 24600
$ws.Range("J61").Value = 24600
$ws.Range("L61").Value = 24600
$ws.Range("N61").Value = -25296
$ws.Range("H68").Value = 42499.5
$ws.Range("J68").Value = 49999
$ws.Range("L68").Value = 49999
$ws.Range("N68").Value = -51497
$ws.Range("H71").Value = 42499.5
$ws.Range("J71").Value = 49999
$ws.Range("L71").Value = 149997
$ws.Range("N71").Value = -157485
$ws.Range("H113").Value = 900
$ws.Range("I113").Value = 800
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 800
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 1370
$ws.Range("N113").Value = -5340
$ws.Range("H134").Value = 5379.7144
$ws.Range("I134").Value = 5379.7144
$ws.Range("K134").Value = 16139.1432
$ws.Range("M134").Value = -13604.1432

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 483.33334
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 483.33334
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 1450.00002
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -2246.00002
$ws.Range("H86").Value = 300
$ws.Range("J86").Value = 400
$ws.Range("L86").Value = 1200
$ws.Range("N86").Value = -3572
$ws.Range("H89").Value = 300
$ws.Range("J89").Value = 400
$ws.Range("L89").Value = 3600
$ws.Range("N89").Value = -15456

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H132").Value = 4139.4443
$ws.Range("I132").Value = 4314.2856
$ws.Range("J132").Value = 3527.5
$ws.Range("K132").Value = 12942.8568
$ws.Range("L132").Value = 10582.5
$ws.Range("M132").Value = -10412.8568
$ws.Range("N132").Value = -15642.5

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2168.375
$ws.Range("J22").Value = 2720.2
$ws.Range("L22").Value = 2720.2
$ws.Range("N22").Value = -3310.2
$ws.Range("H27").Value = 2168.375
$ws.Range("J27").Value = 2720.2
$ws.Range("L27").Value = 2720.2
$ws.Range("N27").Value = -2934.2
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H61").Value = 3167.7693
$ws.Range("I61").Value = 2993.2222
$ws.Range("J61").Value = 3560.5
$ws.Range("K61").Value = 2993.2222
$ws.Range("L61").Value = 3560.5
$ws.Range("M61").Value = -2791.2222
$ws.Range("N61").Value = -3964.5
$ws.Range("H113").Value = 3167.7693
$ws.Range("I113").Value = 2993.2222
$ws.Range("J113").Value = 3560.5
$ws.Range("K113").Value = 2993.2222
$ws.Range("L113").Value = 3560.5
$ws.Range("M113").Value = -823.2222000000002
$ws.Range("N113").Value = -7900.5

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 25650.5
$ws.Range("J81").Value = 1101
$ws.Range("L81").Value = 2202
$ws.Range("N81").Value = -4324
$ws.Range("H84").Value = 25650.5
$ws.Range("J84").Value = 1101
$ws.Range("L84").Value = 11010
$ws.Range("N84").Value = -21618
$ws.Range("H107").Value = 2691
$ws.Range("I107").Value = 1613.75
$ws.Range("K107").Value = 4841.25
$ws.Range("M107").Value = -2921.25
